$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "签署日期" (sign date) header in C1 is being removed from the template;
# clear the cell's content but keep its existing style (s="2").
$ws.Range("C1").ClearContents()

# Move the saved selection/active cell from B9 to C3.
$ws.Range("C3").Select()
